# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.193.48'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.586.82'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'212.06"
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.810.05'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.01"
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.565.80'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").Value = "'64.09"
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '26.210.73'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = "'213.26"
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = "'143.55"
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = "'15.18"
$ws.Range("D30").Value = "'0.0496"
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").Value = '1.339.61'
$ws.Range("E33").Value = '  +4.58%  '
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("D40").Value = "'5.78"
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'0.956"
$ws.Range("E42").Value = '  -13.26%  '
$ws.Range("D43").Value = "'0.769"
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '1.722.33'
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").Value = "'85.70"
$ws.Range("E47").Value = '  -3.52%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'1.48"
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = "'0.0981"
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.0501"
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = '  -0.26%  '
